$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column K (최종점수) values for rows 2-7
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 57.2
$ws.Range("K4").Value = 55.8
$ws.Range("K5").Value = 51.2
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 44.8

# Update column N (MACRO_SCORE) values for rows 2-7
$ws.Range("N2").Value = 85.96878041621773
$ws.Range("N3").Value = 85.96878041621773
$ws.Range("N4").Value = 85.96878041621773
$ws.Range("N5").Value = 85.96878041621773
$ws.Range("N6").Value = 85.96878041621773
$ws.Range("N7").Value = 85.96878041621773
